$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/General-number-as-text storage so values keep their exact
# printed representation (trailing zeros, "%" suffix, etc.) just like the
# source data, instead of Excel auto-converting numeric-looking strings to
# real numbers/percentages.
$cells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "D9",
    "E9",
    "D10",
    "E10",
    "E11",
    "D12",
    "E12",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "E17",
    "D18",
    "E18",
    "E19",
    "E20",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "E25",
    "D26",
    "E26",
    "E27",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "E44",
    "D45",
    "E45",
    "E46",
    "E47",
    "D48",
    "E48",
    "E49",
    "E50",
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin price / volume(1h) figures
$ws.Range("D2").Value = "261.26"
$ws.Range("E2").Value = "1.82%"
$ws.Range("D3").Value = "27.30"
$ws.Range("E3").Value = "1.35%"
$ws.Range("D4").Value = "4.706"
$ws.Range("E4").Value = "1.06%"
$ws.Range("D5").Value = "0.06084"
$ws.Range("E5").Value = "3.19%"
$ws.Range("D6").Value = "6.672"
$ws.Range("E6").Value = "0.98%"
$ws.Range("D7").Value = "0.8462"
$ws.Range("E7").Value = "-0.54%"
$ws.Range("D8").Value = "0.9261"
$ws.Range("D9").Value = "0.1402"
$ws.Range("E9").Value = "1.79%"
$ws.Range("D10").Value = "0.04713"
$ws.Range("E10").Value = "12.20%"
$ws.Range("E11").Value = "1.35%"
$ws.Range("D12").Value = "0.03087"
$ws.Range("E12").Value = "1.14%"
$ws.Range("E13").Value = "-0.36%"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").Value = "0.49%"
$ws.Range("D15").Value = "0.0006085"
$ws.Range("E15").Value = "-94.07%"
$ws.Range("D16").Value = "0.006149"
$ws.Range("E16").Value = "0.81%"
$ws.Range("E17").Value = "-0.60%"
$ws.Range("D18").Value = "3.140"
$ws.Range("E18").Value = "-0.77%"
$ws.Range("E19").Value = "-0.62%"
$ws.Range("E20").Value = "2.24%"
$ws.Range("E21").Value = "-0.53%"
$ws.Range("D22").Value = "4.099"
$ws.Range("E22").Value = "5.11%"
$ws.Range("D23").Value = "0.04246"
$ws.Range("E23").Value = "0.18%"
$ws.Range("D24").Value = "0.001221"
$ws.Range("E24").Value = "0.29%"
$ws.Range("E25").Value = "-8.45%"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").Value = "0.07%"
$ws.Range("E27").Value = "3.45%"
$ws.Range("D40").Value = "0.03873"
$ws.Range("E40").Value = "2.43%"
$ws.Range("D41").Value = "0.1115"
$ws.Range("E41").Value = "1.54%"
$ws.Range("D42").Value = "0.004111"
$ws.Range("E42").Value = "-34.10%"
$ws.Range("D43").Value = "0.01637"
$ws.Range("E43").Value = "15.89%"
$ws.Range("E44").Value = "0.80%"
$ws.Range("D45").Value = "0.00005138"
$ws.Range("E45").Value = "-4.28%"
$ws.Range("E46").Value = "0.07%"
$ws.Range("E47").Value = "19.72%"
$ws.Range("D48").Value = "0.1358"
$ws.Range("E48").Value = "-46.23%"
$ws.Range("E49").Value = "0.07%"
$ws.Range("E50").Value = "0.07%"
